# Updated cryptos list on Wed Jan 31 23:22:33 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($row, $d, $e)
    if ($d -ne $null) {
        $ws.Range("D$row").NumberFormat = "@"
        $ws.Range("D$row").Value = $d
    }
    if ($e -ne $null) {
        $ws.Range("E$row").NumberFormat = "@"
        $ws.Range("E$row").Value = $e
    }
}

Set-Row 2  "42.623.90"  "  -1.43%  "
Set-Row 3  "2.285.14"   "  -3.27%  "
Set-Row 4  "0.999"      "  -0.06%  "
Set-Row 5  "300.51"     "  -2.84%  "
Set-Row 6  "97.02"      "  -6.27%  "
Set-Row 7  "0.503"      "  -1.71%  "
Set-Row 9  $null         "  -3.92%  "
Set-Row 10 "33.45"      "  -6.06%  "
Set-Row 11 $null         "  -2.18%  "
Set-Row 12 "50.26"      "  -5.08%  "
Set-Row 13 $null         "  -0.09%  "
Set-Row 14 "6.66"       "  -3.69%  "
Set-Row 15 "2.637.55"   "  -3.54%  "
Set-Row 16 "15.19"      "  -2.14%  "
Set-Row 17 "2.318.35"   "  -2.08%  "
Set-Row 18 "0.787"      "  -2.73%  "
Set-Row 19 "42.510.77"  "  -1.65%  "
Set-Row 20 $null         "  -1.93%  "
Set-Row 21 "11.47"      "  -3.21%  "
Set-Row 22 $null         "  -4.85%  "
Set-Row 23 "66.75"      "  -1.91%  "
Set-Row 24 "234.73"     "  -2.00%  "
Set-Row 25 $null         "  -4.99%  "
Set-Row 26 $null         "  -4.00%  "
Set-Row 27 $null         "  +0.02%  "
Set-Row 28 $null         "  -4.46%  "
Set-Row 29 "166.21"     "  +3.00%  "
Set-Row 30 "2.06"       "  -11.23%  "
Set-Row 31 "33.74"      "  -7.33%  "
Set-Row 32 $null         "  -3.62%  "
Set-Row 33 "0.999"      "  -0.07%  "
Set-Row 34 $null         "  -4.69%  "
Set-Row 35 $null         "  -3.84%  "
Set-Row 36 "0.0697"     "  -5.30%  "
Set-Row 37 "4.35"       "  -6.73%  "

# Rows 38 and 39 swap: LidoDAOToken <-> Celestia
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-Row 38 "16.24" "  -10.42%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-Row 39 "2.83" "  -7.66%  "

Set-Row 40 $null         "  -7.36%  "
Set-Row 41 $null         "  -4.33%  "
Set-Row 42 $null         "  -2.89%  "
Set-Row 43 "2.48"       "  -4.16%  "
Set-Row 44 "1.961.41"   "  -3.81%  "
Set-Row 46 "17.91"      "  -8.43%  "
Set-Row 47 "9.69"       "  -8.50%  "
Set-Row 48 $null         "  -7.82%  "
Set-Row 49 "53.16"      "  -7.64%  "
Set-Row 50 $null         "  -3.30%  "
Set-Row 51 "2.505.11"   "  -3.52%  "
